$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 5x6 matrix content (A1:F5) since the new data occupies a smaller range
$ws.Range("A1:F5").Clear()

# Write the new 3x4 matrix values
$values = @(
    @(2, 1, -3, -1),
    @(-1, 3, 2, 12),
    @(3, 1, -3, 0)
)

for ($r = 0; $r -lt 3; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

# Update the selected cell to match the target view state
$ws.Range("K6").Select()

$wb.Save()
